$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Table data: No, Tujuan (destination), Satuan, Tarif A, Tarif B, Tarif C, Tarif D
# ---------------------------------------------------------------------------
$data = @(
    @(1,  "Tanah Pinoh Barat", "PP", 650000, 585000, 585000, 585000),
    @(2,  "Tanah Pinoh",       "PP", 400000, 350000, 350000, 350000),
    @(3,  "Sokan",             "PP", 450000, 400000, 400000, 400000),
    @(4,  "Sayan",             "PP", 350000, 300000, 300000, 300000),
    @(5,  "Menukung",          "PP", 550000, 500000, 500000, 500000),
    @(6,  "Ella Hilir",        "PP", 300000, 260000, 260000, 260000),
    @(7,  "Pinoh Selatan",     "PP", 200000, 150000, 150000, 150000),
    @(8,  "Pinoh Utara",       "PP", 200000, 150000, 150000, 150000),
    @(9,  "Belimbing",         "PP", 200000, 150000, 150000, 150000),
    @(10, "Belimbing Hulu",    "PP", 250000, 200000, 200000, 200000)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$lastRow = $startRow + $data.Count - 1

# NO column (A3:A12): centered
$ws.Range("A3:A" + $lastRow).HorizontalAlignment = -4108

# Destination column (B3:B12): left / top / wrap
$colB = $ws.Range("B3:B" + $lastRow)
$colB.HorizontalAlignment = -4131
$colB.VerticalAlignment = -4160
$colB.WrapText = $true

# Satuan column (C3:C12): centered / top / wrap
$colC = $ws.Range("C3:C" + $lastRow)
$colC.HorizontalAlignment = -4108
$colC.VerticalAlignment = -4160
$colC.WrapText = $true

# Tarif columns (D3:G12): number format, right / top / shrink-to-fit, black font
$colDG = $ws.Range("D3:G" + $lastRow)
$colDG.NumberFormat = "#,##0"
$colDG.HorizontalAlignment = -4152
$colDG.VerticalAlignment = -4160
$colDG.ShrinkToFit = $true
$colDG.Font.Color = 0

# Update the active selection to match the finished layout
[void]$ws.Range("D14").Select()

Write-Output "done"
